# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows 10-12: apply the "mtitleStyle" (s=4) look (same as A9) to the row labels ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# --- Updated grading totals ---
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "60/112"

# --- Remove the unused 3rd "Student Ans / Correct Ans" column pair (G:H) ---
$ws.Range("G15:H21").Clear()

# --- Remove the now-unused tail of the 2nd "Student Ans / Correct Ans" pair (D19:E40) ---
$ws.Range("D19:E40").Clear()

# --- Fill in the 2nd pair's remaining Student Ans cells (D16:D18) ---
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("B10").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("C10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option B"

# --- Fill in the Student Ans (column A) for each graded question ---
# s=5 (correctStyle, green) when the student's answer matches, s=6 (incorrectStyle, red) otherwise
function Set-StudentAns($row, $text, $correct) {
    if ($correct) {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $cell = $ws.Range("A" + $row)
    $cell.PasteSpecial(-4122)
    $cell.Value = $text
}

Set-StudentAns 16 "Option A" $true
Set-StudentAns 18 "Option B" $true
Set-StudentAns 19 "Option C" $true
Set-StudentAns 23 "Option D" $true
Set-StudentAns 24 "Option A" $true
Set-StudentAns 26 "Option D" $false
Set-StudentAns 27 "Option A" $true
Set-StudentAns 28 "Option D" $true
Set-StudentAns 29 "Option D" $true
Set-StudentAns 30 "Option B" $true
Set-StudentAns 31 "Option D" $true
Set-StudentAns 32 "Option C" $true
Set-StudentAns 33 "Option A" $false
Set-StudentAns 34 "Option C" $false
Set-StudentAns 37 "Option A" $true
Set-StudentAns 38 "Option A" $true
Set-StudentAns 39 "Option D" $true
